$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Range("B2").Value = 0.06614079833792563
$ws.Range("C2").Value = 0.006265813580010215
$ws.Range("D2").Value = 0.0339659356398244
$ws.Range("E2").Value = 612
$ws.Range("F2").Value = 692
$ws.Range("H2").Value = 1191

# Row 3 (A3 = 1)
$ws.Range("B3").Value = -0.09830159397202665
$ws.Range("C3").Value = 0.007632346204695217
$ws.Range("D3").Value = 0.01206840221748883
$ws.Range("E3").Value = 640
$ws.Range("F3").Value = 664
$ws.Range("H3").Value = 1180

# Row 4 (A4 = 2)
$ws.Range("B4").Value = -0.1262673089474124
$ws.Range("C4").Value = 0.00308434225887888
$ws.Range("D4").Value = 0.01971033593957043
$ws.Range("E4").Value = 629
$ws.Range("F4").Value = 675
$ws.Range("H4").Value = 1205

# Row 5 (A5 = 3)
$ws.Range("B5").Value = -0.03567200542317828
$ws.Range("C5").Value = -0.001020464846548941
$ws.Range("D5").Value = 0.04342196045701304
$ws.Range("E5").Value = 548
$ws.Range("F5").Value = 756

# Row 6 (A6 = 4)
$ws.Range("B6").Value = -0.04709908495982229
$ws.Range("C6").Value = -0.01951975113613413
$ws.Range("D6").Value = 0.03708238541653338
$ws.Range("E6").Value = 639
$ws.Range("F6").Value = 665
$ws.Range("H6").Value = 1529
